# Updated cryptos list on Fri Apr 28 08:28:17 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking price string into a cell as TEXT
# (matches source data where Price/Volume columns are plain inline strings,
# not numbers) while keeping the cell on the default/unstyled format.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '29.432.31'
$ws.Range("E2").Value = '  +1.60%  '

# Row 3
Set-TextValue "D3" '1.917.40'
$ws.Range("E3").Value = '  +1.61%  '

# Row 4
Set-TextValue "D4" '1.002'
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
Set-TextValue "D5" '326.39'
$ws.Range("E5").Value = '  -1.55%  '

# Row 6
$ws.Range("E6").Value = '  +0.14%  '

# Row 7
Set-TextValue "D7" '0.4747'
$ws.Range("E7").Value = '  +2.75%  '

# Row 8
Set-TextValue "D8" '0.4090'
$ws.Range("E8").Value = '  -0.32%  '

# Row 9
Set-TextValue "D9" '47.82'
$ws.Range("E9").Value = '  +1.03%  '

# Row 10
Set-TextValue "D10" '0.08040'
$ws.Range("E10").Value = '  +0.59%  '

# Row 11
Set-TextValue "D11" '1.009'
$ws.Range("E11").Value = '  +1.98%  '

# Row 12
Set-TextValue "D12" '22.56'
$ws.Range("E12").Value = '  +3.89%  '

# Row 13
Set-TextValue "D13" '1.913.74'
$ws.Range("E13").Value = '  -0.22%  '

# Row 14
Set-TextValue "D14" '5.914'
$ws.Range("E14").Value = '  +0.11%  '

# Row 15
Set-TextValue "D15" '7.167'
$ws.Range("E15").Value = '  +1.34%  '

# Row 16
Set-TextValue "D16" '89.94'
$ws.Range("E16").Value = '  +1.10%  '

# Row 17
$ws.Range("E17").Value = '  +0.25%  '

# Row 18
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D18" '0.06607'
$ws.Range("E18").Value = '  +0.67%  '

# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D19" '0.00001033'
$ws.Range("E19").Value = '  +0.42%  '

# Row 20
Set-TextValue "D20" '17.73'
$ws.Range("E20").Value = '  +1.53%  '

# Row 21
Set-TextValue "D21" '1.000'
$ws.Range("E21").Value = '  -0.15%  '

# Row 22
Set-TextValue "D22" '29.436.38'
$ws.Range("E22").Value = '  +1.61%  '

# Row 23
Set-TextValue "D23" '5.525'
$ws.Range("E23").Value = '  +2.24%  '

# Row 24
Set-TextValue "D24" '11.47'
$ws.Range("E24").Value = '  +2.09%  '

# Row 25
Set-TextValue "D25" '2.208'
$ws.Range("E25").Value = '  +0.17%  '

# Row 26
Set-TextValue "D26" '2.138.50'
$ws.Range("E26").Value = '  +0.23%  '

# Row 27
Set-TextValue "D27" '154.91'
$ws.Range("E27").Value = '  -1.53%  '

# Row 28
Set-TextValue "D28" '19.85'
$ws.Range("E28").Value = '  +1.06%  '

# Row 29
Set-TextValue "D29" '6.141'
$ws.Range("E29").Value = '  +13.79%  '

# Row 30
Set-TextValue "D30" '2.129'
$ws.Range("E30").Value = '  +1.26%  '

# Row 31
Set-TextValue "D31" '117.91'
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("E32").Value = '  +10.05%  '

# Row 33
Set-TextValue "D33" '0.09562'
$ws.Range("E33").Value = '  +2.28%  '

# Row 34
Set-TextValue "D34" '1.433'
$ws.Range("E34").Value = '  +1.76%  '

# Row 35
Set-TextValue "D35" '3.559'
$ws.Range("E35").Value = '  -1.15%  '

# Row 36
Set-TextValue "D36" '5.410'
$ws.Range("E36").Value = '  +2.52%  '

# Row 37
Set-TextValue "D37" '0.06100'
$ws.Range("E37").Value = '  +0.77%  '

# Row 38
Set-TextValue "D38" '0.02255'
$ws.Range("E38").Value = '  +1.18%  '

# Row 39
Set-TextValue "D39" '8.303'
$ws.Range("E39").Value = '  +0.49%  '

# Row 40
Set-TextValue "D40" '1.173'
$ws.Range("E40").Value = '  -1.10%  '

# Row 41
Set-TextValue "D41" '0.5898'
$ws.Range("E41").Value = '  +2.14%  '

# Row 42
Set-TextValue "D42" '2.558'
$ws.Range("E42").Value = '  +12.38%  '

# Row 43
Set-TextValue "D43" '0.1843'
$ws.Range("E43").Value = '  +1.32%  '

# Row 44
Set-TextValue "D44" '10.15'
$ws.Range("E44").Value = '  +0.35%  '

# Row 45
Set-TextValue "D45" '0.07928'
$ws.Range("E45").Value = '  +12.92%  '

# Row 46
Set-TextValue "D46" '1.274'
$ws.Range("E46").Value = '  +1.64%  '

# Row 47
Set-TextValue "D47" '0.5558'
$ws.Range("E47").Value = '  +1.75%  '

# Row 48
Set-TextValue "D48" '12.05'
$ws.Range("E48").Value = '  +0.34%  '

# Row 49
Set-TextValue "D49" '1.935'
$ws.Range("E49").Value = '  +1.87%  '

# Row 50
Set-TextValue "D50" '112.93'
$ws.Range("E50").Value = '  +2.01%  '

# Row 51
Set-TextValue "D51" '44.77'
$ws.Range("E51").Value = '  -2.01%  '
